$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 8 (last game row no longer present after the edit)
$ws.Rows.Item(8).Delete()

# Updated data for rows 2-7: D, F, G, H(home_team), I(away_team), J..AF
$data = @(
    @{ Row=2;  D=225.5;  F=112.1651995305164; G=1;  H="Orlando";    I="Washington";   J=0.4891304347826087; K=98.04566705790299; L=113.9056631455399; M=115.7542057902973; N=76.88141627543038; O=0.36016686228482;  P=0.5811090571205006; Q=0.2807246674491393; R=12.59123630672926; S=11.60299295774648; T=0.2137255966353677; U=0.979608729524161;  V=0.9807431626398754; W=10.48471679811085; X=0.4267410015649452; Y=31;   Z=75.34999999999999; AA=0.4883040935672515; AB=0.4998633419594793; AC=-5.45;  AD=0.1024090272268591;  AE=0.09658417421970092; AF=0.4639018941231757 }
    @{ Row=3;  D=236;    F=114.1485719874804; G=14; H="Atlanta";    I="Detroit";      J=0.528169014084507;  K=99.36740023474178; L=114.0723200312989; M=118.0094581377152; N=75.26235328638495; O=0.3518038928012519; P=0.5693385172143978; Q=0.2686073943661972; R=11.74988262910798; S=11.8726917057903;  T=0.2165450899843505; U=0.9969307597159863; V=0.9574346620970108; W=11.35669312351145; X=0.3575899843505477; Y=38;   Z=75.5;               AA=0.3821428571428572; AB=0.5007783244838535; AC=-2.9;   AD=0.2868701557620145;  AE=0.2253176010150409;  AF=0.513206101967522  }
    @{ Row=4;  D=218.5;  F=112.7513988037816; G=3;  H="Brooklyn";   I="Cleveland";    J=0.4680039138943248; K=96.25072351919738; L=116.3805807447424; M=113.2076596565695; N=75.7050935751495;  O=0.3840830600038587; P=0.5974148176731624; Q=0.2640336677599847; R=11.91571483696701; S=12.66660235384912; T=0.2107690526721976; U=0.9847283738321537; V=0.9626442643961158; W=10.82799689195475; X=0.5828670654061354; Y=48.5; Z=76.84999999999999; AA=0.4719251336898396; AB=0.5015132535795598; AC=1.43;   AD=0.03986769961597369; AE=0.03242224178403757; AF=0.4779323265780384}
    @{ Row=5;  D=232.5;  F=113.3591549295775; G=13; H="NewOrleans"; I="SanAntonio";   J=0.5363339457032048; K=99.52394366197183; L=113.0718309859155; M=117.3830985915493; N=76.10140845070423; O=0.346161971830986;  P=0.569549295774648;  Q=0.2578239436619718; R=12.72887323943662; S=12.2556338028169;  T=0.2052394366197183; U=0.9900362875945631; V=1.045650279378181;  W=10.75281481949624; X=0.3732394366197183; Y=34;   Z=75.30000000000001; AA=0.4494949494949495; AB=0.5008197623188704; AC=-0.41;  AD=0.2309774288349089;  AE=0.1797412949921753;  AF=0.512552038787907  }
    @{ Row=6;  D=239.5;  F=119.3405712050078; G=5;  H="Sacramento"; I="Boston";       J=0.5142857142857142; K=98.91913145539905; L=119.2578051643193; M=115.1836952269171; N=78.2709996087637;  O=0.4495860719874802; P=0.6070291471048512; Q=0.2723900625978091; R=11.76524843505477; S=11.37936228482003; T=0.207889622456964;  U=1.042275731048103;  V=0.9776736366503012; W=12.46603764290631; X=0.6430946791862284; Y=44.5; Z=77.25;              AA=0.5727272727272728; AB=0.4805047624224144; AC=0.8799999999999999; AD=0.009573168958799755; AE=0.01630063657407408; AF=0.4725093050704769}
    @{ Row=7;  D=238.5;  F=115.3462441314554; G=7;  H="LAClippers"; I="OklahomaCity"; J=0.5179968701095461; K=98.86306729264476; L=115.5945618153365; M=114.7527875586855; N=75.04909037558686; O=0.3771229460093897; P=0.5790099765258215; Q=0.2690244522691706; R=11.58657081377152; S=12.68953442879499; T=0.2143623337245696; U=1.007390778440658;  V=1.033706398810785;  W=13.40429650336272; X=0.5103677621283256; Y=38;   Z=75.40000000000001; AA=0.3808823529411764; AB=0.4880580593205012; AC=-3.905; AD=0.06129130655821043; AE=0.04027874228395063; AF=0.5258694122900343}
)

foreach ($r in $data) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
    $ws.Range("U$row").Value = $r.U
    $ws.Range("V$row").Value = $r.V
    $ws.Range("W$row").Value = $r.W
    $ws.Range("X$row").Value = $r.X
    $ws.Range("Y$row").Value = $r.Y
    $ws.Range("Z$row").Value = $r.Z
    $ws.Range("AA$row").Value = $r.AA
    $ws.Range("AB$row").Value = $r.AB
    $ws.Range("AC$row").Value = $r.AC
    $ws.Range("AD$row").Value = $r.AD
    $ws.Range("AE$row").Value = $r.AE
    $ws.Range("AF$row").Value = $r.AF
}
